$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Fgf7"
$ws.Range("C2").Value = "Fgfr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 18.32032066666667
$ws.Range("H2").Value = 54.960962
$ws.Range("I2").Value = 0.9393883922888365
$ws.Range("J2").Value = 0.9570930138987986
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.495057333333333
$ws.Range("N2").Value = 16.485172
$ws.Range("O2").Value = 0.8161989011161211
$ws.Range("P2").Value = 0.8403205285996808
$ws.Range("Q2").Value = 100.6712124283849
$ws.Range("R2").Value = 906.040911855464
$ws.Range("S2").Value = 0.7667277735073881
$ws.Range("T2").Value = 0.8042649073585001

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Fgf7"
$ws.Range("C3").Value = "Fgfr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 18.32032066666667
$ws.Range("H3").Value = 54.960962
$ws.Range("I3").Value = 0.9393883922888365
$ws.Range("J3").Value = 0.9570930138987986
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.657666
$ws.Range("N3").Value = 1.972998
$ws.Range("O3").Value = 0.09768528951377062
$ws.Range("P3").Value = 0.1005722428790014
$ws.Range("Q3").Value = 12.048652011564
$ws.Range("R3").Value = 108.437868104076
$ws.Range("S3").Value = 0.09176442706661052
$ws.Range("T3").Value = 0.09625699105162543

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Fgf7"
$ws.Range("C4").Value = "Fgfr3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 18.32032066666667
$ws.Range("H4").Value = 54.960962
$ws.Range("I4").Value = 0.9393883922888365
$ws.Range("J4").Value = 0.9570930138987986
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.5797745
$ws.Range("N4").Value = 1.159549
$ws.Range("O4").Value = 0.08611580937010824
$ws.Range("P4").Value = 0.0591072285213179
$ws.Range("Q4").Value = 10.62165475435633
$ws.Range("R4").Value = 63.729928526138
$ws.Range("S4").Value = 0.0808961917148379
$ws.Range("T4").Value = 0.05657111548867318

# Row 5
$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "Fgf7"
$ws.Range("C5").Value = "Fgfr3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.09978633333333332
$ws.Range("H5").Value = 0.299359
$ws.Range("I5").Value = 0.005116620224500324
$ws.Range("J5").Value = 0.005213052994737072
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.495057333333333
$ws.Range("N5").Value = 16.485172
$ws.Range("O5").Value = 0.8161989011161211
$ws.Range("P5").Value = 0.8403205285996808
$ws.Range("Q5").Value = 0.5483316227497776
$ws.Range("R5").Value = 4.934984604747999
$ws.Range("S5").Value = 0.004176179804665686
$ws.Range("T5").Value = 0.004380635448155605

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Fgf7"
$ws.Range("C6").Value = "Fgfr3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.09978633333333332
$ws.Range("H6").Value = 0.299359
$ws.Range("I6").Value = 0.005116620224500324
$ws.Range("J6").Value = 0.005213052994737072
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.657666
$ws.Range("N6").Value = 1.972998
$ws.Range("O6").Value = 0.09768528951377062
$ws.Range("P6").Value = 0.1005722428790014
$ws.Range("Q6").Value = 0.065626078698
$ws.Range("R6").Value = 0.5906347082819999
$ws.Range("S6").Value = 0.0004998185279623282
$ws.Range("T6").Value = 0.0005242884319278024

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Fgf7"
$ws.Range("C7").Value = "Fgfr3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.09978633333333332
$ws.Range("H7").Value = 0.299359
$ws.Range("I7").Value = 0.005116620224500324
$ws.Range("J7").Value = 0.005213052994737072
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.5797745
$ws.Range("N7").Value = 1.159549
$ws.Range("O7").Value = 0.08611580937010824
$ws.Range("P7").Value = 0.0591072285213179
$ws.Range("Q7").Value = 0.05785357151516666
$ws.Range("R7").Value = 0.347121429091
$ws.Range("S7").Value = 0.0004406218918723104
$ws.Range("T7").Value = 0.0003081291146536648

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Fgf7"
$ws.Range("C8").Value = "Fgfr3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.082285
$ws.Range("H8").Value = 2.16457
$ws.Range("I8").Value = 0.05549498748666317
$ws.Range("J8").Value = 0.03769393310646423
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 5.495057333333333
$ws.Range("N8").Value = 16.485172
$ws.Range("O8").Value = 0.8161989011161211
$ws.Range("P8").Value = 0.8403205285996808
$ws.Range("Q8").Value = 5.947218126006666
$ws.Range("R8").Value = 35.68330875603999
$ws.Range("S8").Value = 0.04529494780406737
$ws.Range("T8").Value = 0.03167498579302503

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Fgf7"
$ws.Range("C9").Value = "Fgfr3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.082285
$ws.Range("H9").Value = 2.16457
$ws.Range("I9").Value = 0.05549498748666317
$ws.Range("J9").Value = 0.03769393310646423
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.657666
$ws.Range("N9").Value = 1.972998
$ws.Range("O9").Value = 0.09768528951377062
$ws.Range("P9").Value = 0.1005722428790014
$ws.Range("Q9").Value = 0.7117820468099999
$ws.Range("R9").Value = 4.27069228086
$ws.Range("S9").Value = 0.00542104391919777
$ws.Range("T9").Value = 0.003790963395448152

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Fgf7"
$ws.Range("C10").Value = "Fgfr3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.082285
$ws.Range("H10").Value = 2.16457
$ws.Range("I10").Value = 0.05549498748666317
$ws.Range("J10").Value = 0.03769393310646423
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.5797745
$ws.Range("N10").Value = 1.159549
$ws.Range("O10").Value = 0.08611580937010824
$ws.Range("P10").Value = 0.0591072285213179
$ws.Range("Q10").Value = 0.6274812447324999
$ws.Range("R10").Value = 2.50992497893
$ws.Range("S10").Value = 0.004778995763398027
$ws.Range("T10").Value = 0.002227983917991052
